$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: B2 is a plain, default-styled text cell (no custom NumberFormat).
# We use its .Style as a template to strip off any "@" text NumberFormat
# we temporarily apply below, so cell formatting stays byte-identical to
# the original (only the cell contents change).
$textStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "26.196.13"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "1.675.16"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  -0.78%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.88"
$ws.Range("D5").Style = $textStyle
$ws.Range("E5").Value = "  -3.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5266"
$ws.Range("D6").Style = $textStyle
$ws.Range("E6").Value = "  -3.96%  "
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2652"
$ws.Range("D8").Style = $textStyle
$ws.Range("E8").Value = "  -3.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06288"
$ws.Range("D9").Style = $textStyle
$ws.Range("E9").Value = "  -2.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.37"
$ws.Range("D10").Style = $textStyle
$ws.Range("E10").Value = "  -3.31%  "
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "1.660.11"
$ws.Range("E12").Value = "  -2.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.455"
$ws.Range("D13").Style = $textStyle
$ws.Range("E13").Value = "  -2.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5615"
$ws.Range("D14").Style = $textStyle
$ws.Range("E14").Value = "  -4.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.97"
$ws.Range("D15").Style = $textStyle
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000007996"
$ws.Range("D16").Style = $textStyle
$ws.Range("E16").Value = "  -4.85%  "
$ws.Range("D17").Value = "25.959.65"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.812"
$ws.Range("D19").Style = $textStyle
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "187.91"
$ws.Range("D20").Style = $textStyle
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.43"
$ws.Range("D21").Style = $textStyle
$ws.Range("E21").Value = "  -5.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.180"
$ws.Range("D22").Style = $textStyle
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("D23").Style = $textStyle
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.76"
$ws.Range("D24").Style = $textStyle
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1247"
$ws.Range("D25").Style = $textStyle
$ws.Range("E25").Value = "  -5.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.564"
$ws.Range("D26").Style = $textStyle
$ws.Range("E26").Value = "  -4.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.02"
$ws.Range("D27").Style = $textStyle
$ws.Range("E27").Value = "  +1.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06182"
$ws.Range("D28").Style = $textStyle
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.360"
$ws.Range("D29").Style = $textStyle
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.285"
$ws.Range("D30").Style = $textStyle
$ws.Range("E30").Value = "  -3.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.491"
$ws.Range("D31").Style = $textStyle
$ws.Range("E31").Value = "  -3.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.437"
$ws.Range("D32").Style = $textStyle
$ws.Range("E32").Value = "  -4.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.630"
$ws.Range("D33").Style = $textStyle
$ws.Range("E33").Value = "  -3.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.001"
$ws.Range("D34").Style = $textStyle
$ws.Range("E34").Value = "  -4.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6065"
$ws.Range("D35").Style = $textStyle
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.406"
$ws.Range("D36").Style = $textStyle
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.737"
$ws.Range("D37").Style = $textStyle
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.087"
$ws.Range("D38").Style = $textStyle
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01612"
$ws.Range("D39").Style = $textStyle
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("D40").Value = "1.081.95"
$ws.Range("E40").Value = "  -3.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8718"
$ws.Range("D41").Style = $textStyle
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("E42").Value = "  -1.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.93"
$ws.Range("D43").Style = $textStyle
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("D44").Value = "1.824.56"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000107"
$ws.Range("D45").Style = $textStyle
$ws.Range("E45").Value = "  -1.63%  "
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").Style = $textStyle
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.013"
$ws.Range("D48").Style = $textStyle
$ws.Range("E48").Value = "  -3.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05236"
$ws.Range("D49").Style = $textStyle
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4257"
$ws.Range("D50").Style = $textStyle
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.969"
$ws.Range("D51").Style = $textStyle
$ws.Range("E51").Value = "  -3.00%  "
